$wb = $excel.ActiveWorkbook

# Map of current sheet name -> new sheet name
$renames = @{
    "GP1" = "GP01"
    "GP2" = "GP02"
    "BP1" = "BP01"
    "BP2" = "BP02"
    "BP3" = "BP03"
    "BP4" = "BP04"
    "BP5" = "BP05"
    "BP6" = "BP06"
    "BP7" = "BP07"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $newName = $renames[$oldName]

        # Update the KPI title text in cell A1, e.g.
        # "KPI GP1 - Global Perf 1" -> "KPI GP01 - Global Perf 1"
        $cellValue = $ws.Range("A1").Value2
        if ($cellValue -ne $null) {
            $updatedValue = $cellValue -replace ("KPI " + $oldName + " "), ("KPI " + $newName + " ")
            $ws.Range("A1").Value2 = $updatedValue
        }

        $ws.Name = $newName
    }
}
